$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 header values
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 values (B2 cleared, C2:E2 updated)
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 5.1090777001051748
$ws.Range("D2").Value = 0.15909593397877569
$ws.Range("E2").Value = 3.2166496700074267

# Row 3 values
$ws.Range("B3").Value = 0.22614458846323032
$ws.Range("C3").Value = 6.5387440478270467
$ws.Range("D3").Value = 0.29079876538191207
$ws.Range("E3").Value = 4.4255690177005844

# Update selection to match new selected range
$ws.Range("B1:E3").Select() | Out-Null
